# Commit: "rename sheet to match bean"
#
# The worksheet (and the workbook's hidden "ExternalData_2" defined name that
# points at it) is renamed from the old ad-hoc name to the name of the
# Salesforce bean/Step it feeds:
#   "SEAWARE-BOOKINGS-Promos"  ->  "Step00-Upsert-PromoItem"
#
# Renaming via the Worksheet.Name property automatically keeps any
# definedName/formula references that point at the sheet (e.g.
# ExternalData_2) in sync.

$wb = $excel.ActiveWorkbook

$oldName = "SEAWARE-BOOKINGS-Promos"
$newName = "Step00-Upsert-PromoItem"

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq $oldName) {
        $ws = $sheet
        break
    }
}

if ($ws -eq $null) {
    # Fall back to whichever sheet is active if the expected name isn't found.
    $ws = $wb.ActiveSheet
}

$ws.Name = $newName
